$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.342.23'
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').Value = '3.378.82'
$ws.Range('E3').Value = '  -2.85%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').Value = '''141.52'
$ws.Range('E6').Value = '  -4.90%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.376.82'
$ws.Range('E8').Value = '  -2.87%  '
$ws.Range('E9').Value = '  -3.17%  '
$ws.Range('E10').Value = '  +4.76%  '
$ws.Range('D11').Value = '''0.133'
$ws.Range('E11').Value = '  -6.95%  '
$ws.Range('D12').Value = '''0.404'
$ws.Range('E12').Value = '  -5.07%  '
$ws.Range('D13').Value = '3.953.17'
$ws.Range('E13').Value = '  -2.76%  '
$ws.Range('E14').Value = '  -7.41%  '
$ws.Range('D15').Value = '''29.43'
$ws.Range('E15').Value = '  -7.33%  '
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '65.371.45'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '3.379.97'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = '''10.29'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = '''6.08'
$ws.Range('E20').Value = '  -6.26%  '
$ws.Range('E21').Value = '  -5.89%  '
$ws.Range('D22').Value = '''412.25'
$ws.Range('E22').Value = '  -6.39%  '
$ws.Range('D23').Value = '''0.577'
$ws.Range('E23').Value = '  -5.97%  '
$ws.Range('D24').Value = '''76.82'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '3.522.09'
$ws.Range('E26').Value = '  -2.55%  '
$ws.Range('E27').Value = '  -10.10%  '
$ws.Range('D28').Value = '''9.18'
$ws.Range('E28').Value = '  -6.30%  '
$ws.Range('D29').Value = '''7.76'
$ws.Range('E29').Value = '  -7.89%  '
$ws.Range('E30').Value = '  -3.57%  '
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '''0.160'
$ws.Range('E32').Value = '  -5.27%  '
$ws.Range('E33').Value = '  -9.01%  '
$ws.Range('D34').Value = '''24.28'
$ws.Range('E34').Value = '  -4.64%  '
$ws.Range('D35').Value = '3.375.72'
$ws.Range('E35').Value = '  -2.58%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').Value = '''5.51'
$ws.Range('E37').Value = '  -9.32%  '
$ws.Range('E38').Value = '  -7.77%  '
$ws.Range('D39').Value = '''0.999'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').Value = '''7.47'
$ws.Range('E40').Value = '  -5.97%  '
$ws.Range('D41').Value = '''168.31'
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('E42').Value = '  -4.72%  '
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('E44').Value = '  -7.89%  '
$ws.Range('E45').Value = '  -11.16%  '
$ws.Range('D46').Value = '''45.29'
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '''26.36'
$ws.Range('E47').Value = '  -9.94%  '
$ws.Range('E48').Value = '  -5.85%  '
$ws.Range('E49').Value = '  -6.39%  '
$ws.Range('D50').Value = '''2.25'
$ws.Range('E50').Value = '  -8.71%  '
$ws.Range('D51').Value = '''0.914'
$ws.Range('E51').Value = '  -7.38%  '
